# "changes in ms report 09/07/25"
# Clears the report footer row (page number, confidentiality notice and
# run-date stamp) from the "Report Design" sheet while leaving the
# surrounding cell formatting untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 footer text:
#   A20  -> "Page 1 of 1"
#   D20:H20 (merged) -> "Confidential Information of Fresenius Kabi.  Do not copy or distribute."
#   K20  -> "Run Date: Jul 07, 2025"
$ws.Range("A20").ClearContents()
$ws.Range("D20:H20").ClearContents()
$ws.Range("K20").ClearContents()

# Leave the cursor where the author last left it.
$ws.Range("J23").Select() | Out-Null
